$d = $word.ActiveDocument

function Replace-FirstText($findText, $newText) {
    # Locate the first occurrence of $findText (case-sensitive, whole string,
    # not a wildcard pattern) starting from the top of the document, then
    # overwrite just that run's text. Setting .Text directly (instead of
    # passing $newText as Find.Execute's Replace argument) avoids Word's
    # "replace" autocorrect pass that would otherwise turn straight quotes
    # into curly ones.
    $r = $d.Content
    $found = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if ($found) {
        $r.Text = $newText
    } else {
        Write-Output "NOT FOUND: $findText"
    }
}

# 1. Equipment needed sentence
Replace-FirstText "Set for each group of students: a glass or a cup, a plate, one stick, soap, water, coffee, black pepper. Both the plates and the water should be clean." "Fixé pour chaque groupe d'étudiants : un verre ou une tasse, une assiette, un bâton, savon, eau, café, poivre noir. Les plaques et l'eau devraient être propres."

# 2. "None" -> "Aucun"
Replace-FirstText "None" "Aucun"

# 3. "Video d'introduction" -> "Vidéo d'introduction"
Replace-FirstText "Video d'introduction" "Vidéo d'introduction"

# 4. First "VIDEO PAUSE" / "Experiment:" / "Filling the glasses with water" block
Replace-FirstText "VIDEO PAUSE" "PAUSE VIDÉO"
Replace-FirstText "Experiment:" "Expérience :"
Replace-FirstText "Filling the glasses with water" "Remplir les verres d'eau"

# 5. "Faciliter le processus, susciter des pensées" -> "Aider le processus, provoquer des réflexions" (first occurrence)
Replace-FirstText "Faciliter le processus, susciter des pensées" "Aider le processus, provoquer des réflexions"

# 6. "Fill the glasses or cups up to the very top" -> "Remplissez les verres ou les tasses jusqu'au bord"
Replace-FirstText "Fill the glasses or cups up to the very top" "Remplissez les verres ou les tasses jusqu'au bord"

# 7. "Introduction of the second experiment" -> "Introduction de la deuxième expérience"
Replace-FirstText "Introduction of the second experiment" "Introduction de la deuxième expérience"
